$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 391.8
$ws.Range("J41").Value = 1166.3334
$ws.Range("L41").Value = 1166.3334
$ws.Range("N41").Value = -2046.3334
$ws.Range("H107").Value = 1752.2632
$ws.Range("I107").Value = 1514.9375
$ws.Range("J107").Value = 3018
$ws.Range("K107").Value = 1514.9375
$ws.Range("L107").Value = 3018
$ws.Range("M107").Value = 405.0625
$ws.Range("N107").Value = -6858
$ws.Range("H111").Value = 2132
$ws.Range("I111").Value = 2498
$ws.Range("K111").Value = 7494
$ws.Range("M111").Value = -4427
$ws.Range("H132").Value = 2801.4517
$ws.Range("I132").Value = 2563.2856
$ws.Range("K132").Value = 7689.8568
$ws.Range("M132").Value = -5159.8568
$ws.Range("H134").Value = 58945.918
$ws.Range("J134").Value = 58945.918
$ws.Range("L134").Value = 58945.918
$ws.Range("N134").Value = -69085.91800000001
$ws.Range("H135").Value = 3167.5
$ws.Range("I135").Value = 3167.5
$ws.Range("K135").Value = 28507.5
$ws.Range("M135").Value = -25972.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1450.8462
$ws.Range("I110").Value = 1710.875
$ws.Range("K110").Value = 1710.875
$ws.Range("M110").Value = 334.125
$ws.Range("H122").Value = 1854.0769
$ws.Range("I122").Value = 1794.8889
$ws.Range("J122").Value = 1987.25
$ws.Range("K122").Value = 5384.6667
$ws.Range("L122").Value = 5961.75
$ws.Range("M122").Value = -2934.6667
$ws.Range("N122").Value = -10861.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9999.6
$ws.Range("I16").Value = 9999.75
$ws.Range("K16").Value = 9999.75
$ws.Range("M16").Value = -9712.75
$ws.Range("H31").Value = 1243.1316
$ws.Range("I31").Value = 797.56525
$ws.Range("J31").Value = 1926.3334
$ws.Range("K31").Value = 797.56525
$ws.Range("L31").Value = 1926.3334
$ws.Range("M31").Value = -502.56525
$ws.Range("N31").Value = -2516.3334
$ws.Range("H34").Value = 1243.1316
$ws.Range("I34").Value = 797.56525
$ws.Range("J34").Value = 1926.3334
$ws.Range("K34").Value = 797.56525
$ws.Range("L34").Value = 1926.3334
$ws.Range("M34").Value = -595.56525
$ws.Range("N34").Value = -2330.3334
$ws.Range("H62").Value = 303135.3
$ws.Range("I62").Value = 203499.4
$ws.Range("J62").Value = 402771.2
$ws.Range("K62").Value = 203499.4
$ws.Range("L62").Value = 402771.2
$ws.Range("M62").Value = -202875.4
$ws.Range("N62").Value = -404019.2
$ws.Range("H65").Value = 303135.3
$ws.Range("I65").Value = 203499.4
$ws.Range("J65").Value = 402771.2
$ws.Range("K65").Value = 1017497
$ws.Range("L65").Value = 2013856
$ws.Range("M65").Value = -1014377
$ws.Range("N65").Value = -2020096
$ws.Range("H80").Value = 89332.664
$ws.Range("J80").Value = 89332.664
$ws.Range("L80").Value = 89332.664
$ws.Range("N80").Value = -91578.664
$ws.Range("H83").Value = 89332.664
$ws.Range("J83").Value = 89332.664
$ws.Range("L83").Value = 267997.992
$ws.Range("N83").Value = -279229.992
$ws.Range("H113").Value = 9999.6
$ws.Range("I113").Value = 9999.75
$ws.Range("K113").Value = 9999.75
$ws.Range("M113").Value = -7829.75
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("H134").Value = 21899.166
$ws.Range("I134").Value = 13799.667
$ws.Range("K134").Value = 41399.001
$ws.Range("M134").Value = -38864.001
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2794.1667
$ws.Range("I5").Value = 887
$ws.Range("J5").Value = 3175.6
$ws.Range("K5").Value = 2661
$ws.Range("L5").Value = 9526.799999999999
$ws.Range("M5").Value = -2549
$ws.Range("N5").Value = -9750.799999999999
$ws.Range("H107").Value = 948.7
$ws.Range("I107").Value = 824.2727
$ws.Range("J107").Value = 1100.7778
$ws.Range("K107").Value = 2472.8181
$ws.Range("L107").Value = 3302.3334
$ws.Range("M107").Value = -552.8181
$ws.Range("N107").Value = -7142.3334
$ws.Range("H114").Value = 2008.5
$ws.Range("I114").Value = 663
$ws.Range("J114").Value = 2176.6875
$ws.Range("K114").Value = 1989
$ws.Range("L114").Value = 6530.0625
$ws.Range("M114").Value = 1265
$ws.Range("N114").Value = -13038.0625
$ws.Range("H135").Value = 2794.1667
$ws.Range("I135").Value = 887
$ws.Range("J135").Value = 3175.6
$ws.Range("K135").Value = 7983
$ws.Range("L135").Value = 28580.4
$ws.Range("M135").Value = -5448
$ws.Range("N135").Value = -33650.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15638.75
$ws.Range("J15").Value = 15638.75
$ws.Range("L15").Value = 15638.75
$ws.Range("N15").Value = -16214.75
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("H80").Value = 11672.2
$ws.Range("I80").Value = 7351.6
$ws.Range("J80").Value = 15992.8
$ws.Range("K80").Value = 7351.6
$ws.Range("L80").Value = 15992.8
$ws.Range("M80").Value = -6353.6
$ws.Range("N80").Value = -17988.8
$ws.Range("H81").Value = 15638.75
$ws.Range("J81").Value = 15638.75
$ws.Range("L81").Value = 15638.75
$ws.Range("N81").Value = -17634.75
$ws.Range("H83").Value = 11672.2
$ws.Range("I83").Value = 7351.6
$ws.Range("J83").Value = 15992.8
$ws.Range("K83").Value = 36758
$ws.Range("L83").Value = 79964
$ws.Range("M83").Value = -31766
$ws.Range("N83").Value = -89948
$ws.Range("H84").Value = 15638.75
$ws.Range("J84").Value = 15638.75
$ws.Range("L84").Value = 46916.25
$ws.Range("N84").Value = -56900.25
$ws.Range("H113").Value = 2503.2307
$ws.Range("I113").Value = 2071.7
$ws.Range("K113").Value = 2071.7
$ws.Range("M113").Value = 98.30000000000018
$ws.Range("M60").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 633.25
$ws.Range("I16").Value = 679.1429000000001
$ws.Range("J16").Value = 312
$ws.Range("K16").Value = 679.1429000000001
$ws.Range("L16").Value = 312
$ws.Range("M16").Value = -509.1429000000001
$ws.Range("N16").Value = -652
$ws.Range("H22").Value = 1431.7
$ws.Range("I22").Value = 2239.8
$ws.Range("J22").Value = 1162.3334
$ws.Range("K22").Value = 2239.8
$ws.Range("L22").Value = 1162.3334
$ws.Range("M22").Value = -1944.8
$ws.Range("N22").Value = -1752.3334
$ws.Range("H27").Value = 1431.7
$ws.Range("I27").Value = 2239.8
$ws.Range("J27").Value = 1162.3334
$ws.Range("K27").Value = 2239.8
$ws.Range("L27").Value = 1162.3334
$ws.Range("M27").Value = -2132.8
$ws.Range("N27").Value = -1376.3334
$ws.Range("H136").Value = 4838.6
$ws.Range("I136").Value = 4509.5454
$ws.Range("J136").Value = 5743.5
$ws.Range("K136").Value = 13528.6362
$ws.Range("L136").Value = 17230.5
$ws.Range("M136").Value = -10978.6362
$ws.Range("N136").Value = -22330.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3925.8
$ws.Range("I62").Value = 3740.8333
$ws.Range("J62").Value = 4665.6665
$ws.Range("K62").Value = 3740.8333
$ws.Range("L62").Value = 4665.6665
$ws.Range("M62").Value = -3116.8333
$ws.Range("N62").Value = -5913.6665
$ws.Range("H65").Value = 3925.8
$ws.Range("I65").Value = 3740.8333
$ws.Range("J65").Value = 4665.6665
$ws.Range("K65").Value = 18704.1665
$ws.Range("L65").Value = 23328.3325
$ws.Range("M65").Value = -15584.1665
$ws.Range("N65").Value = -29568.3325
$ws.Range("H75").Value = 32236.166
$ws.Range("I75").Value = 44949
$ws.Range("K75").Value = 44949
$ws.Range("M75").Value = -44013
$ws.Range("H78").Value = 32236.166
$ws.Range("I78").Value = 44949
$ws.Range("K78").Value = 134847
$ws.Range("M78").Value = -130167
$ws.Range("H132").Value = 6081.324
$ws.Range("I132").Value = 3943.8635
$ws.Range("K132").Value = 11831.5905
$ws.Range("M132").Value = -9301.5905

Write-Output "Applied Twintania_Profits market price updates"